$d = $word.ActiveDocument

function Set-CellText($table, $row, $col, $text) {
    $cell = $table.Rows.Item($row).Cells.Item($col)
    $cell.Range.Text = $text
}

# Table 1 - header info
$t1 = $d.Tables.Item(1)
Set-CellText $t1 2 4 "05/12/2021, 01:46 PM"
Set-CellText $t1 3 4 "TAC S1 E4"

# Table 2 - Capacity
$t2 = $d.Tables.Item(2)
Set-CellText $t2 2 4 "484"

# Table 3 - Chilled Water Circuit
$t3 = $d.Tables.Item(3)
Set-CellText $t3 2 4 "292"
Set-CellText $t3 3 4 "12"
Set-CellText $t3 4 4 "7"
Set-CellText $t3 5 4 "1+1"
Set-CellText $t3 6 4 "3.4"
Set-CellText $t3 7 4 "200"
Set-CellText $t3 10 4 "standard"

# Table 4 - Cooling Water Circuit
$t4 = $d.Tables.Item(4)
Set-CellText $t4 2 4 "484"
Set-CellText $t4 3 4 "32"
Set-CellText $t4 4 4 "38.8"
Set-CellText $t4 5 4 "1+1/1"
Set-CellText $t4 7 4 "9.1"
Set-CellText $t4 8 4 "250"
Set-CellText $t4 11 4 "standard"

# Table 5 - Steam Circuit
$t5 = $d.Tables.Item(5)
Set-CellText $t5 3 4 "3238.8"
Set-CellText $t5 6 4 "200"
Set-CellText $t5 7 4 "40"

# Table 6 - Electrical Data
$t6 = $d.Tables.Item(6)
Set-CellText $t6 2 4 "415 V( ±10%), 50 Hz (±5%), 3 Phase+N"
Set-CellText $t6 3 4 "13.4"
Set-CellText $t6 4 4 "5.5( 14 )"
Set-CellText $t6 5 4 "0.3( 1.4 )"
# Remove MOP (row 7) and MCA (row 8) rows
$t6.Rows.Item(7).Delete()
$t6.Rows.Item(7).Delete()

# Table 7 - Physical Data
$t7 = $d.Tables.Item(7)
Set-CellText $t7 2 4 "4910"
Set-CellText $t7 3 4 "2305"
Set-CellText $t7 4 4 "3395"
Set-CellText $t7 5 4 "13.7"
Set-CellText $t7 6 4 "12.7"
Set-CellText $t7 7 4 "21.8"
Set-CellText $t7 8 4 "10"
Set-CellText $t7 9 4 "4690"

# Footnotes: remove note 1 and note 6 paragraphs entirely, renumber the remaining items 2-5 to 1-4
# NOTE: use $d.Content.Paragraphs (not the bare $d.Paragraphs collection) for
# indexed access here -- after navigating into Tables.Item(...) above, this
# interop layer's bare $d.Paragraphs.Item(n) indexer can get anchored to the
# last-touched table range; going through $d.Content.Paragraphs keeps indices
# relative to the whole document as expected.

# Renumber notes 2-5 -> 1-4 in place (keeps the existing run/formatting, just
# swaps the text) before removing notes 1 and 6 entirely.
# A paragraph Range.Text includes the trailing paragraph-mark ("`r"), so
# compare against both the bare and "`r"-suffixed forms to be safe.
for ($i = 1; $i -le $d.Content.Paragraphs.Count; $i++) {
    $p = $d.Content.Paragraphs.Item($i)
    $txt = $p.Range.Text
    if ($txt -eq "2. This selection is valid for insulated chiller only." -or $txt -eq "2. This selection is valid for insulated chiller only.`r") {
        $p.Range.Text = "1. This selection is valid for insulated chiller only."
    } elseif ($txt -eq "3. For non-insulated chiller, the Capacity and Heat source consumption will vary." -or $txt -eq "3. For non-insulated chiller, the Capacity and Heat source consumption will vary.`r") {
        $p.Range.Text = "2. For non-insulated chiller, the Capacity and Heat source consumption will vary."
    } elseif ($txt -eq "4. Plant Room Temperature should be from +5 deg C to +45 deg C" -or $txt -eq "4. Plant Room Temperature should be from +5 deg C to +45 deg C`r") {
        $p.Range.Text = "3. Plant Room Temperature should be from +5 deg C to +45 deg C"
    } elseif ($txt -eq "5. Please contact Thermax representative / Office for customised specifications." -or $txt -eq "5. Please contact Thermax representative / Office for customised specifications.`r") {
        $p.Range.Text = "4. Please contact Thermax representative / Office for customised specifications."
    }
}

# Delete whole paragraphs (including the paragraph mark) from the end first so
# earlier paragraph indices are not invalidated.
for ($i = $d.Content.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Content.Paragraphs.Item($i)
    $txt = $p.Range.Text
    if ($txt -eq "6. Try Reducing Cooling water flow`r" -or $txt -eq "6. Try Reducing Cooling water flow") {
        $p.Range.Delete()
    }
}
for ($i = $d.Content.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Content.Paragraphs.Item($i)
    $txt = $p.Range.Text
    if ($txt -eq "1. This is an ARI selection`r" -or $txt -eq "1. This is an ARI selection") {
        $p.Range.Delete()
    }
}
